# update w/ harris poll (11/24)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the survey "week" value (12 -> 13) for the four existing Harris
#     rows (111-114) that belong to the same wave as the new poll below ---
$ws.Range("C111").Value = 13
$ws.Range("C112").Value = 13
$ws.Range("C113").Value = 13
$ws.Range("C114").Value = 13

# --- Append three new Harris Interactive poll rows (121-123) ---

# Row 121
$ws.Cells.Item(121, 1).Value = 34        # id
$ws.Cells.Item(121, 2).Value = 2021      # year
$ws.Cells.Item(121, 3).Value = 13        # week
$ws.Cells.Item(121, 4).Value = 11        # month
$ws.Cells.Item(121, 5).Value = 21        # day
$ws.Cells.Item(121, 6).Value = "harris"     # firm
$ws.Cells.Item(121, 7).Value = "online"     # collectmode
$ws.Cells.Item(121, 8).Value = "included"   # unsure
$ws.Cells.Item(121, 9).Value = 2120       # n
$ws.Cells.Item(121, 10).Value = 1         # c_poutou
$ws.Cells.Item(121, 11).Value = 1         # c_arthaud
$ws.Cells.Item(121, 12).Value = 10        # c_melenchon
$ws.Cells.Item(121, 13).Value = 2         # c_roussel
$ws.Cells.Item(121, 14).Value = 1         # c_montebourg
$ws.Cells.Item(121, 15).Value = 8         # c_jadot
$ws.Cells.Item(121, 16).Value = 5         # c_hidalgo
$ws.Cells.Item(121, 17).Value = 23        # c_macron
$ws.Cells.Item(121, 20).Value = 14        # c_bertrand
$ws.Cells.Item(121, 21).Value = "T_0.5"   # c_lassalle
$ws.Cells.Item(121, 22).Value = 2         # c_daignant
$ws.Cells.Item(121, 23).Value = 16        # c_lepen
$ws.Cells.Item(121, 24).Value = 16        # c_zemmour
$ws.Cells.Item(121, 25).Value = "T_0.5"   # c_asselineau
$ws.Cells.Item(121, 27).Value = 1         # c_philippot

# Row 122
$ws.Cells.Item(122, 1).Value = 34
$ws.Cells.Item(122, 2).Value = 2021
$ws.Cells.Item(122, 3).Value = 13
$ws.Cells.Item(122, 4).Value = 11
$ws.Cells.Item(122, 5).Value = 21
$ws.Cells.Item(122, 6).Value = "harris"
$ws.Cells.Item(122, 7).Value = "online"
$ws.Cells.Item(122, 8).Value = "included"
$ws.Cells.Item(122, 9).Value = 2121
$ws.Cells.Item(122, 10).Value = 1
$ws.Cells.Item(122, 11).Value = 1
$ws.Cells.Item(122, 12).Value = 10
$ws.Cells.Item(122, 13).Value = 2
$ws.Cells.Item(122, 14).Value = 1
$ws.Cells.Item(122, 15).Value = 9
$ws.Cells.Item(122, 16).Value = 5
$ws.Cells.Item(122, 17).Value = 24
$ws.Cells.Item(122, 18).Value = 11        # c_pecresse
$ws.Cells.Item(122, 21).Value = "T_0.5"
$ws.Cells.Item(122, 22).Value = 2
$ws.Cells.Item(122, 23).Value = 16
$ws.Cells.Item(122, 24).Value = 16
$ws.Cells.Item(122, 25).Value = "T_0.5"
$ws.Cells.Item(122, 27).Value = 2

# Row 123
$ws.Cells.Item(123, 1).Value = 34
$ws.Cells.Item(123, 2).Value = 2021
$ws.Cells.Item(123, 3).Value = 13
$ws.Cells.Item(123, 4).Value = 11
$ws.Cells.Item(123, 5).Value = 21
$ws.Cells.Item(123, 6).Value = "harris"
$ws.Cells.Item(123, 7).Value = "online"
$ws.Cells.Item(123, 8).Value = "included"
$ws.Cells.Item(123, 9).Value = 2122
$ws.Cells.Item(123, 10).Value = 1
$ws.Cells.Item(123, 11).Value = 1
$ws.Cells.Item(123, 12).Value = 10
$ws.Cells.Item(123, 13).Value = 2
$ws.Cells.Item(123, 14).Value = 1
$ws.Cells.Item(123, 15).Value = 9
$ws.Cells.Item(123, 16).Value = 5
$ws.Cells.Item(123, 17).Value = 24
$ws.Cells.Item(123, 19).Value = 10        # c_barnier
$ws.Cells.Item(123, 21).Value = "T_0.5"
$ws.Cells.Item(123, 22).Value = 2
$ws.Cells.Item(123, 23).Value = 16
$ws.Cells.Item(123, 24).Value = 17
$ws.Cells.Item(123, 25).Value = "T_0.5"
$ws.Cells.Item(123, 27).Value = 2

# Move the selection to reflect where the user ended up after entering
# the new data (bottom-right of the newly extended used range).
$ws.Range("AB123").Select()
